$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> (Type1, Type2) data, matching the order pokemon appear in rows 2..52
$types = @{
    2  = @('Planta','Hielo')
    3  = @('Fantasma','Acero')
    4  = @('Roca','Acero')
    5  = @('Roca','Hiero')
    6  = @('Lucha','Normal')
    7  = @('Lucha','Fuego')
    8  = @('Acero','Phiquico')
    9  = @('Fuego','Fantasma')
    10 = @('Volador','Fuego')
    11 = @('Agua','Hielo')
    12 = @('Dragon','Volador')
    13 = @('Electrico')
    14 = @('Siniestro')
    15 = @('Agua')
    16 = @('Dragon','Tierra')
    17 = @('Hielo','Fantasma')
    18 = @('Lucha','Phiquico')
    19 = @('Bicho','Electrico')
    20 = @('Dragon','Tierra')
    21 = @('Hada','Phiquico')
    22 = @('Veneno','Fantasma')
    23 = @('Hielo')
    24 = @('Tierra','Volador')
    25 = @('Planta')
    26 = @('Electrico')
    27 = @('Acero','Lucha')
    28 = @('Agua')
    29 = @('Electrico')
    30 = @('Roca')
    31 = @('Acero','Electrico')
    32 = @('Acero','Phiquico')
    33 = @('Agua')
    34 = @('Hada','Fantasma')
    35 = @('Fantasma')
    36 = @('Normal')
    37 = @('Electrico','Phiquico')
    39 = @('Veneno','Fuego')
    40 = @('Planta')
    41 = @('Bicho','Veneno')
    42 = @('Planta')
    43 = @('Siniestro','Fantasma')
    44 = @('Hada')
    45 = @('Hada','Volador')
    46 = @('Fuego')
    47 = @('Roca','Siniestro')
    48 = @('Siniestro')
    49 = @('Agua')
    50 = @('Planta','Veneno')
    51 = @('Normal','Hada')
    52 = @('Bicho','Volador')
}

foreach ($row in $types.Keys) {
    $vals = $types[$row]
    if ($vals.Count -ge 1) {
        $ws.Range("K$row").Value = $vals[0]
    }
    if ($vals.Count -ge 2) {
        $ws.Range("L$row").Value = $vals[1]
    }
}

# Row 38 (Reuniclus) only has a second type populated (Phiquico), no first type
$ws.Range("L38").Value = 'Phiquico'

# Update the selection to match the author's final cursor position
$ws.Range("M19").Select()
